$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the header formatting used by the rest of row 1 (bold, bordered,
# centered) by copying the format from the adjacent header cell.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# New header cells for the team record columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record for every player row (2-52)
$ws.Range("AD2:AD52").Value = 86
$ws.Range("AE2:AE52").Value = 76
$ws.Range("AF2:AF52").Value = 0
